$d = $word.ActiveDocument

function Find-ParagraphIndex($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. Insert two new bullet paragraphs (ListParagraph style, ilvl=3, numId=1)
#    right after "Emails numerical code to each member..." and before
#    "Search a buy-order".
# ---------------------------------------------------------------------------
$emailsIdx = Find-ParagraphIndex("Emails numerical code to each member in a group and buy-order creator has master list to check with")
$emailsPara = $d.Paragraphs.Item($emailsIdx)

$ins1 = $emailsPara.Range.Duplicate
$ins1.Collapse(0)
$ins1.InsertParagraphAfter()

$barcodePara = $d.Paragraphs.Item($emailsIdx + 1)
$barcodePara.Range.Text = "(Potential) add a serializable barcode to determine each product uniquely"

$ins2 = $barcodePara.Range.Duplicate
$ins2.Collapse(0)
$ins2.InsertParagraphAfter()

$supportPara = $d.Paragraphs.Item($emailsIdx + 2)
$supportPara.Range.Text = "Support number for drivers"

# ---------------------------------------------------------------------------
# 2. Move the hidden "_GoBack" bookmark from the end of "Backend updates the
#    a list of carts..." paragraph to the end of the new "Support number for
#    drivers" paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# NOTE: Bookmarks.Add with a collapsed range sitting exactly on the last
# character position of a paragraph (Range.End - 1) is mishandled by this
# COM shim (it silently resets to the very start of the document). Work
# around it by bookmarking a temporary placeholder character, then clearing
# the bookmarked range's text -- this leaves a proper zero-length bookmark
# in the correct spot, right before the paragraph mark.
$bmHost = $supportPara.Range.Duplicate
$bmHost.Collapse(0)
$bmHost.MoveEnd(1, -1)
$bmHost.InsertAfter("X")

$d.Bookmarks.Add("_GoBack", $bmHost)
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Range.Text = ""

# ---------------------------------------------------------------------------
# 3. Move the <w:lastRenderedPageBreak/> marker from the run beginning
#    "(optional) Location, request Location" to the run beginning "Picture ".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("(optional) Location, request Location", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Location, request Location", 2) | Out-Null

Write-Output "done"
